# Atualização de bases das ligas, do dia: 02-05-2024 às 20:28
#
# The two fixtures stored in rows 53 and 55 of the "Montenegro Prva Liga"
# sheet had been mixed up: the match id / teams / odds data that belongs to
# fixture 7279987 (FK Jezero vs FK Arsenal) was sitting in row 53 while the
# data for fixture 6815334 (Sutjeska Niksic vs FK Mornar Bar) was sitting in
# row 55. This swaps the two rows' data back into the correct rows, leaving
# the row-index column (A) and the unchanged columns (C, D, H, V) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Montenegro Prva Liga")

# Columns that differ between the two rows, in sheet order.
$cols = @("B","E","F","G","I","J","K","L","M","N","O","P","Q","R","S","T","U","W","X","Y","Z","AA","AB")

# Target values for row 53 (fixture 7279987: FK Jezero vs FK Arsenal).
$row53 = @{
    "B" = 7279987; "E" = "FK Jezero"; "F" = "FK Arsenal"; "G" = 1; "I" = "D"
    "J" = 2.1; "K" = 3; "L" = 3.25; "M" = 2.05; "N" = 3; "O" = 3.4
    "P" = -0.25; "Q" = 1.8; "R" = 2; "S" = 2; "T" = 1.925; "U" = 1.875
    "W" = 2; "X" = -1; "Y" = -0.5; "Z" = 0.5; "AA" = 0; "AB" = 0
}

# Target values for row 55 (fixture 6815334: Sutjeska Niksic vs FK Mornar Bar).
$row55 = @{
    "B" = 6815334; "E" = "Sutjeska Niksic"; "F" = "FK Mornar Bar"; "G" = 0; "I" = "A"
    "J" = 1.444; "K" = 4; "L" = 6.5; "M" = 1.444; "N" = 4; "O" = 6.5
    "P" = -1.25; "Q" = 2; "R" = 1.8; "S" = 2.5; "T" = 2; "U" = 1.8
    "W" = -1; "X" = 5.5; "Y" = -1; "Z" = 0.8; "AA" = -1; "AB" = 0.8
}

foreach ($c in $cols) {
    $ws.Range($c + "53").Value = $row53[$c]
    $ws.Range($c + "55").Value = $row55[$c]
}
